# Reworked section 3. Added open set figure.
#
# This script reproduces the reproducible parts of the target edit using the
# PowerPoint COM/automation object model:
#   1. Recolor the "baseline" series line in the Section-3 FROC chart
#      (slide 1 / chart1.xml) from green (00B050) to purple (813FB7).
#   2. Nudge/resize that same chart's graphic frame a couple of EMUs to the
#      left and a little wider (matches the new off/ext on slide 1).
#   3. Bump the cached "today" footer date field (datetimeFigureOut) from
#      05.09.2017 to 06.09.2017 everywhere it is cached - the slide master
#      and all 11 slide layouts.
#
# (Internal chart axis-id bookkeeping (<c:axId>/<c:crossAx>) and shape
# book-keeping ids/modIds (<p:cNvPr id>, <p14:modId>) are not settable
# through the PowerPoint object model - PowerPoint itself assigns/refreshes
# those transparently - so they are intentionally left alone here.)

$p = $ppt.ActivePresentation

# --- 1 & 2: chart on slide 1 -------------------------------------------
$slide1 = $p.Slides.Item(1)
$chartShape = $slide1.Shapes.Item(1)
$chart1 = $chartShape.Chart

# Series 7 is "baseline, public (0.001)" - the last <c:ser> in chart1.xml,
# the one whose line color changes in the diff.
$series = $chart1.SeriesCollection().Item(7)
$series.Format.Line.ForeColor.RGB = 12009345   # RGB(0x81,0x3F,0xB7) == 813FB7

# Resize/reposition the chart's graphic frame: off -137,0 / ext 2880000x2330451 (EMU)
#                                            -> off -3947,0 / ext 2887620x2330451 (EMU)
# TextFrame/Shape geometry is expressed in points (1 pt = 12700 EMU) over COM.
$chartShape.Left = -3947 / 12700.0
$chartShape.Width = 2887620 / 12700.0

# --- 3: footer date field bump (05.09.2017 -> 06.09.2017) ---------------
function Set-DatumText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Datumsplatzhalter*") {
            $shp.TextFrame.TextRange.Text = "06.09.2017"
        }
    }
}

Set-DatumText $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Set-DatumText $layout.Shapes
}
